# Phase 4.4: Fix remaining validation errors
# Update the "Experimental" flag to false and refresh the "Date" timestamp
# on the Metadata sheet of the CodeSystem workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# B7 holds the "Experimental" value. Enter it with a leading apostrophe so
# the engine stores it as literal text "false" (a shared string) rather
# than coercing it to a native boolean, then restore the original cell
# formatting (borrowed from the neighboring "Status" value cell) so the
# quote-prefix styling introduced by the apostrophe entry is discarded.
$expCell = $ws.Range("B7")
$expCell.Value = "'false"
$ws.Range("B6").Copy()
$expCell.PasteSpecial(-4122)

# B8 holds the "Date" value; simply refresh the timestamp text.
$ws.Range("B8").Value = "2025-10-03T16:37:46+01:00"
